# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to the Leve profit-tracking tables
# (columns H:N) across the ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets,
# as produced by the scheduled market-data refresh runner.

$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 126
$ws.Range("H126").Value = 35000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 35000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 35000
$ws.Range("N126").Value = -44880
# Row 129
$ws.Range("H129").Value = 1326.2
$ws.Range("I129").Value = 292.25
$ws.Range("J129").Value = 1584.6875
$ws.Range("K129").Value = 876.75
$ws.Range("L129").Value = 4754.0625
$ws.Range("M129").Value = 4123.25
$ws.Range("N129").Value = -14754.0625

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 23418.645
$ws.Range("I32").Value = 4635.8433
$ws.Range("J32").Value = 110502.55
$ws.Range("K32").Value = 4635.8433
$ws.Range("L32").Value = 110502.55
$ws.Range("M32").Value = -4348.8433
$ws.Range("N32").Value = -111076.55
# Row 74
$ws.Range("H74").Value = 3805.8696
$ws.Range("I74").Value = 1039.2
$ws.Range("J74").Value = 8993.375
$ws.Range("K74").Value = 1039.2
$ws.Range("L74").Value = 8993.375
$ws.Range("M74").Value = -165.2
$ws.Range("N74").Value = -10741.375
# Row 77
$ws.Range("H77").Value = 3805.8696
$ws.Range("I77").Value = 1039.2
$ws.Range("J77").Value = 8993.375
$ws.Range("K77").Value = 5196
$ws.Range("L77").Value = 44966.875
$ws.Range("M77").Value = -828
$ws.Range("N77").Value = -53702.875
# Row 102
$ws.Range("H102").Value = 2041.1428
$ws.Range("I102").Value = 1777.6
$ws.Range("J102").Value = 2700
$ws.Range("K102").Value = 1777.6
$ws.Range("L102").Value = 2700
$ws.Range("M102").Value = -155.5999999999999
$ws.Range("N102").Value = -5944
# Row 109
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
# Row 122
$ws.Range("H122").Value = 55000
$ws.Range("I122").Value = 55000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 165000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -162550
$ws.Range("N122").ClearContents()

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# Row 40
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
# Row 86
$ws.Range("H86").Value = 1590.7273
$ws.Range("I86").Value = 1642.5714
$ws.Range("J86").Value = 1500
$ws.Range("K86").Value = 1642.5714
$ws.Range("L86").Value = 1500
$ws.Range("M86").Value = -519.5714
$ws.Range("N86").Value = -3746
# Row 89
$ws.Range("H89").Value = 1590.7273
$ws.Range("I89").Value = 1642.5714
$ws.Range("J89").Value = 1500
$ws.Range("K89").Value = 8212.857
$ws.Range("L89").Value = 7500
$ws.Range("M89").Value = -2596.857
$ws.Range("N89").Value = -18732
# Row 96
$ws.Range("H96").Value = 23985.6
$ws.Range("I96").Value = 16642.666
$ws.Range("J96").Value = 35000
$ws.Range("K96").Value = 16642.666
$ws.Range("L96").Value = 35000
$ws.Range("M96").Value = -13896.666
$ws.Range("N96").Value = -40492
# Row 99
$ws.Range("H99").Value = 1644.5883
$ws.Range("I99").Value = 1381.4615
$ws.Range("J99").Value = 2499.75
$ws.Range("K99").Value = 1381.4615
$ws.Range("L99").Value = 2499.75
$ws.Range("M99").Value = 116.5385000000001
$ws.Range("N99").Value = -5495.75
# Row 107
$ws.Range("H107").Value = 2637.5
$ws.Range("I107").Value = 2562.4
$ws.Range("J107").Value = 3013
$ws.Range("K107").Value = 2562.4
$ws.Range("L107").Value = 3013
$ws.Range("M107").Value = -642.4000000000001
$ws.Range("N107").Value = -6853
# Row 115
$ws.Range("H115").Value = 35000
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 35000
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 35000
$ws.Range("N115").Value = -38134

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1632.8182
$ws.Range("I16").Value = 1468.5
$ws.Range("J16").Value = 1830
$ws.Range("K16").Value = 1468.5
$ws.Range("L16").Value = 1830
$ws.Range("M16").Value = -1181.5
$ws.Range("N16").Value = -2404
# Row 105
$ws.Range("H105").Value = 958.25
$ws.Range("I105").Value = 1015.7143
$ws.Range("J105").Value = 877.8
$ws.Range("K105").Value = 1015.7143
$ws.Range("L105").Value = 877.8
$ws.Range("M105").Value = 731.2857
$ws.Range("N105").Value = -4371.8
# Row 113
$ws.Range("H113").Value = 1632.8182
$ws.Range("I113").Value = 1468.5
$ws.Range("J113").Value = 1830
$ws.Range("K113").Value = 1468.5
$ws.Range("L113").Value = 1830
$ws.Range("M113").Value = 701.5
$ws.Range("N113").Value = -6170
# Row 132
$ws.Range("H132").Value = 1793.4242
$ws.Range("I132").Value = 1248.3572
$ws.Range("J132").Value = 4845.8
$ws.Range("K132").Value = 3745.0716
$ws.Range("L132").Value = 14537.4
$ws.Range("M132").Value = -1215.0716
$ws.Range("N132").Value = -19597.4
# Row 134
$ws.Range("H134").Value = 2300.2327
$ws.Range("I134").Value = 1355.6389
$ws.Range("J134").Value = 7158.143
$ws.Range("K134").Value = 4066.9167
$ws.Range("L134").Value = 21474.429
$ws.Range("M134").Value = -1531.9167
$ws.Range("N134").Value = -26544.429

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 1849.75
$ws.Range("I113").Value = 1699.5
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1699.5
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 470.5
$ws.Range("N113").Value = -6340
# Row 122
$ws.Range("H122").Value = 1094.5
$ws.Range("I122").Value = 800
$ws.Range("J122").Value = 1683.5
$ws.Range("K122").Value = 2400
$ws.Range("L122").Value = 5050.5
$ws.Range("M122").Value = 50
$ws.Range("N122").Value = -9950.5
# Row 126
$ws.Range("H126").Value = 2908.5
$ws.Range("I126").Value = 2782.1667
$ws.Range("J126").Value = 2955.875
$ws.Range("K126").Value = 8346.500100000001
$ws.Range("L126").Value = 8867.625
$ws.Range("M126").Value = -5876.500100000001
$ws.Range("N126").Value = -13807.625
# Row 132
$ws.Range("H132").Value = 3502.1143
$ws.Range("I132").Value = 2565.8
$ws.Range("J132").Value = 5842.9
$ws.Range("K132").Value = 7697.400000000001
$ws.Range("L132").Value = 17528.7
$ws.Range("M132").Value = -5167.400000000001
$ws.Range("N132").Value = -22588.7

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 880.0714
$ws.Range("I46").Value = 603.1429000000001
$ws.Range("J46").Value = 1157
$ws.Range("K46").Value = 603.1429000000001
$ws.Range("L46").Value = 1157
$ws.Range("M46").Value = -415.1429000000001
$ws.Range("N46").Value = -1533
# Row 61
$ws.Range("H61").Value = 1168.75
$ws.Range("I61").Value = 950
$ws.Range("J61").Value = 1533.3334
$ws.Range("K61").Value = 950
$ws.Range("L61").Value = 1533.3334
$ws.Range("M61").Value = -748
$ws.Range("N61").Value = -1937.3334
# Row 92
$ws.Range("H92").Value = 29428.572
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 29428.572
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 29428.572
$ws.Range("N92").Value = -34420.572
# Row 100
$ws.Range("H100").Value = 2455.9443
$ws.Range("I100").Value = 1486.4286
$ws.Range("J100").Value = 3072.9092
$ws.Range("K100").Value = 1486.4286
$ws.Range("L100").Value = 3072.9092
$ws.Range("M100").Value = -945.4286
$ws.Range("N100").Value = -4154.9092
# Row 113
$ws.Range("H113").Value = 1168.75
$ws.Range("I113").Value = 950
$ws.Range("J113").Value = 1533.3334
$ws.Range("K113").Value = 950
$ws.Range("L113").Value = 1533.3334
$ws.Range("M113").Value = 1220
$ws.Range("N113").Value = -5873.3334
# Row 132
$ws.Range("H132").Value = 3080
$ws.Range("I132").Value = 2034.5312
$ws.Range("J132").Value = 7859.2856
$ws.Range("K132").Value = 6103.5936
$ws.Range("L132").Value = 23577.8568
$ws.Range("M132").Value = -3573.5936
$ws.Range("N132").Value = -28637.8568
# Row 136
$ws.Range("H136").Value = 4864.2383
$ws.Range("I136").Value = 2099.7058
$ws.Range("J136").Value = 16613.5
$ws.Range("K136").Value = 6299.117400000001
$ws.Range("L136").Value = 49840.5
$ws.Range("M136").Value = -3749.117400000001
$ws.Range("N136").Value = -54940.5

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 3210.861
$ws.Range("I81").Value = 2027.1666
$ws.Range("J81").Value = 4394.5557
$ws.Range("K81").Value = 4054.3332
$ws.Range("L81").Value = 8789.1114
$ws.Range("M81").Value = -2993.3332
$ws.Range("N81").Value = -10911.1114
# Row 84
$ws.Range("H84").Value = 3210.861
$ws.Range("I84").Value = 2027.1666
$ws.Range("J84").Value = 4394.5557
$ws.Range("K84").Value = 20271.666
$ws.Range("L84").Value = 43945.557
$ws.Range("M84").Value = -14967.666
$ws.Range("N84").Value = -54553.557
# Row 100
$ws.Range("H100").Value = 827
$ws.Range("I100").Value = 897.25
$ws.Range("J100").Value = 733.3333
$ws.Range("K100").Value = 1794.5
$ws.Range("L100").Value = 1466.6666
$ws.Range("M100").Value = -1253.5
$ws.Range("N100").Value = -2548.6666
# Row 113
$ws.Range("H113").Value = 457.10526
$ws.Range("I113").Value = 483.3125
$ws.Range("J113").Value = 317.33334
$ws.Range("K113").Value = 1449.9375
$ws.Range("L113").Value = 952.0000200000001
$ws.Range("M113").Value = 720.0625
$ws.Range("N113").Value = -5292.00002
# Row 126
$ws.Range("H126").Value = 114133.336
$ws.Range("I126").Value = 128275
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 384825
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -382355
$ws.Range("N126").Value = -7940
# Row 132
$ws.Range("H132").Value = 3524.639
$ws.Range("I132").Value = 3442.889
$ws.Range("J132").Value = 3769.889
$ws.Range("K132").Value = 10328.667
$ws.Range("L132").Value = 11309.667
$ws.Range("M132").Value = -7798.667000000001
$ws.Range("N132").Value = -16369.667
# Row 136
$ws.Range("H136").Value = 3100.9575
$ws.Range("I136").Value = 1226.9117
$ws.Range("J136").Value = 8002.3076
$ws.Range("K136").Value = 3680.7351
$ws.Range("L136").Value = 24006.9228
$ws.Range("M136").Value = -1130.7351
$ws.Range("N136").Value = -29106.9228

